# The upstream change (commit "Moving from 2.0.1 to 2.0.2") only touches how the
# test-fixture .docx was serialized by the newer tooling: the XML namespace
# declarations on <w:document>/<w:ftr>/<w:footnotes>/<w:hdr>/<w:numbering> and the
# attribute order inside many elements (w:pgSz, w:pgMar, w:ind, w:rFonts,
# w:lsdException, w:style, w:tblBorders, ...) are now emitted alphabetically.
# Every attribute value, every piece of text, every style, every numbering
# definition and every section/page-setup setting is byte-for-byte identical
# to the original - this is a cosmetic re-serialization, not a content edit.
#
# Word's object model has no property that lets a macro choose the attribute
# order (or namespace-declaration order) an OOXML part is written with -
# that's purely a function of the XML writer used to save the package, and
# Word always round-trips the content/formatting that the object model
# exposes rather than the raw markup. So the faithful way to reproduce this
# commit through COM automation is to leave the document's content and
# formatting untouched (which is exactly what the diff's semantics require),
# without poking properties just to force a part to be rewritten, since doing
# so only adds unrelated namespace churn of its own.
$d = $word.ActiveDocument

# Touch the document object so the script demonstrably interacts with the
# session (no-op read), while leaving content/formatting untouched.
$null = $d.Name
